$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.046.27"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.563.22"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.20"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.65"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.60"
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.45"
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  -4.43%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.957.06"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "2.539.98"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.843"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "43.086.39"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.84"
$ws.Range("E19").Value = "  +3.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.57"
$ws.Range("E20").Value = "  -3.74%  "
$ws.Range("D21").Value = "0.0₃0961"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.44"
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.72"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.07"
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.80"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.92"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.21"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  -4.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "154.52"
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.41"
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0805"
$ws.Range("E35").Value = "  +1.88%  "
$ws.Range("E36").Value = "  +2.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.04"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  +5.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.118"
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.39"
$ws.Range("E41").Value = "  -5.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.95"
$ws.Range("E42").Value = "  +3.54%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("E45").Value = "  -1.76%  "
$ws.Range("D46").Value = "2.001.96"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "83.27"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("D49").Value = "2.812.34"
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.36"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.92"
$ws.Range("E51").Value = "  +0.90%  "
